$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = -1.9770267617657384
$ws.Range("C3").Value = -1.9497280051420887

$ws.Range("B4").Value = -2.1247586523022473
$ws.Range("C4").Value = -2.0719191754432713

$ws.Range("B5").Value = -2.2515035901787703
$ws.Range("C5").Value = -2.1801205335574969

$ws.Range("B6").Value = -2.2194752216690619
$ws.Range("C6").Value = -2.1290762581173475

$ws.Range("B7").Value = -2.2626980658252043
$ws.Range("C7").Value = -2.1559522509384146

$ws.Range("B9").Value = 0.35605176398309912
$ws.Range("C9").Value = 0.35825921679285777

$ws.Range("B13").Value = 0.061789755953912097
$ws.Range("C13").Value = 0.062282729873806488

$ws.Range("B14").Value = 1.0792832465898308
$ws.Range("C14").Value = 1.0782819005691589

$ws.Range("C16").Value = -0.147226934973433

$ws.Range("C17").Value = -0.021307115571185753

$ws.Range("C18").Value = 0.11776015566474256

$ws.Range("C19").Value = -1.5403318388201941

$ws.Range("C20").Value = 0.03712047798613282

$ws.Range("B21").Value = 8259
$ws.Range("C21").Value = 10222
